# Null-ordering auto-detection: the two rows that represent the "missing
# from right" and "missing from left" outer-join cases for key=3 (rows 5/6)
# and for key=5 (rows 9/10) were written in the wrong relative order.
# Swap each pair of rows back into the order the database would return.
#
# xlPasteValues = -4163. Using Copy + PasteSpecial(xlPasteValues) moves the
# cell's stored value/type intact (round-trips numeric-looking text such as
# "12" as text rather than re-parsing it as a number the way a plain
# `.Value = "12"` assignment would) while leaving the destination cell's
# existing style index untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("multicolsWithNulls")

function Swap-Range($rangeA, $rangeB, $tempRange) {
    $ws.Range($rangeA).Copy()
    $ws.Range($tempRange).PasteSpecial(-4163)

    $ws.Range($rangeB).Copy()
    $ws.Range($rangeA).PasteSpecial(-4163)

    $ws.Range($tempRange).Copy()
    $ws.Range($rangeB).PasteSpecial(-4163)

    $ws.Range($tempRange).Clear()
}

# Row 5 <-> Row 6, columns B:I (key=3 pair)
Swap-Range "B5:I5" "B6:I6" "B1000:I1000"

# Row 9 <-> Row 10, columns C:I (key=5 pair; column B already matches in both)
Swap-Range "C9:I9" "C10:I10" "C1000:I1000"
